$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 165; existing rows 165-238 shift down to 166-239.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new record.
$ws.Cells.Item(165, 1).Value2  = 3
$ws.Cells.Item(165, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(165, 3).Value2  = "Coquimbo"
$ws.Cells.Item(165, 4).Value2  = 44609
$ws.Cells.Item(165, 5).Value2  = 5
$ws.Cells.Item(165, 6).Value2  = 100112001
$ws.Cells.Item(165, 7).Value2  = "Berenjena"
$ws.Cells.Item(165, 8).Value2  = "Sin especificar"
$ws.Cells.Item(165, 9).Value2  = "Primera"
$ws.Cells.Item(165, 10).Value2 = 110
$ws.Cells.Item(165, 11).Value2 = 9500
$ws.Cells.Item(165, 12).Value2 = 10000
$ws.Cells.Item(165, 13).Value2 = 9727
$ws.Cells.Item(165, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(165, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(165, 16).Value2 = 162
$ws.Cells.Item(165, 17).Value2 = 60
$ws.Cells.Item(165, 18).Value2 = "Hortaliza"
